$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 425, shifting existing rows 425:529 down to 426:530
$ws.Range("A425").EntireRow.Insert()

# Populate the newly inserted row 425 with the new data record
$ws.Range("A425").Value2 = 10
$ws.Range("B425").Value2 = "Vega Modelo de Temuco"
$ws.Range("C425").Value2 = "La Araucanía"
$ws.Range("D425").Value2 = 44932
$ws.Range("E425").Value2 = 9
$ws.Range("F425").Value2 = 100112008
$ws.Range("G425").Value2 = "Coliflor"
$ws.Range("H425").Value2 = "Sin especificar"
$ws.Range("I425").Value2 = "Primera"
$ws.Range("J425").Value2 = 500
$ws.Range("K425").Value2 = 1200
$ws.Range("L425").Value2 = 1200
$ws.Range("M425").Value2 = 1200
$ws.Range("N425").Value2 = "$/unidad"
$ws.Range("O425").Value2 = "Provincia de Cautín"
$ws.Range("P425").Value2 = 1200
$ws.Range("Q425").Value2 = 1
$ws.Range("R425").Value2 = "Hortaliza"
